# Adicion de columna Toxicidad a diccionario
# Adds a new row (43) to the data dictionary sheet describing the new
# "toxicidad" field.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 43

$ws.Cells.Item($newRow, 1).Value = "toxicidad"
$ws.Cells.Item($newRow, 2).Value = 50
$ws.Cells.Item($newRow, 3).Value = "string"
$ws.Cells.Item($newRow, 4).Value = "Nivel de toxicidad de pla planta."

# Match the author's final view/selection state: scrolled so row 10 is at
# the top, with D43 (the last cell they typed into) selected.
$ws.Range("D43").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
